# Update marksheet totals: correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - total right-answer points
$ws.Range("B11").Value = 5

# "Total" row - total score and "correct/total" display string
$ws.Range("B12").Value = 95
$ws.Range("E12").Value = "95/140"
